$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.423.10'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.010.19'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'251.97"
$ws.Range('E5').Value = '  +2.99%  '
$ws.Range('E6').Value = '  -3.03%  '
$ws.Range('D7').Value = "'61.63"
$ws.Range('E7').Value = '  +12.94%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.371"
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('D10').Value = "'58.62"
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = "'0.104"
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('D13').Value = "'0.897"
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = "'14.83"
$ws.Range('E14').Value = '  +3.89%  '
$ws.Range('D15').Value = '2.302.33'
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('D16').Value = "'20.25"
$ws.Range('E16').Value = '  +15.61%  '
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '2.026.04'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '36.373.04'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = "'71.96"
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = '0.0₃0863'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('D22').Value = "'5.27"
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = "'234.43"
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('D24').Value = "'2.76"
$ws.Range('E24').Value = '  +21.33%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').Value = "'9.57"
$ws.Range('E27').Value = '  +3.45%  '
$ws.Range('D28').Value = "'163.77"
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').Value = "'19.63"
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('D32').Value = "'0.110"
$ws.Range('E32').Value = '  +22.22%  '
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').Value = "'4.60"
$ws.Range('E34').Value = '  +5.76%  '
$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('E36').Value = '  +10.75%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = "'5.91"
$ws.Range('E39').Value = '  +17.55%  '
$ws.Range('E40').Value = '  +13.85%  '
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('D42').Value = "'2.78"
$ws.Range('E42').Value = '  +23.18%  '
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0216"
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'8.01"
$ws.Range('E46').Value = '  +7.54%  '
$ws.Range('D47').Value = "'16.73"
$ws.Range('E47').Value = '  +7.71%  '
$ws.Range('D48').Value = "'94.33"
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('D49').Value = '1.425.55'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('D51').Value = "'46.93"
$ws.Range('E51').Value = '  +2.64%  '
